$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 250.75
$ws.Range("I12").Value = 183.66667
$ws.Range("J12").Value = 452
$ws.Range("K12").Value = 183.66667
$ws.Range("L12").Value = 452
$ws.Range("M12").Value = -13.66667000000001
$ws.Range("N12").Value = -792
$ws.Range("H18").Value = 3333965
$ws.Range("I18").Value = 5000200
$ws.Range("K18").Value = 5000200
$ws.Range("M18").Value = -4999916
$ws.Range("H33").Value = 875.3333
$ws.Range("I33").Value = 950.4
$ws.Range("K33").Value = 950.4
$ws.Range("M33").Value = -721.4
$ws.Range("H43").Value = 9498.571
$ws.Range("J43").Value = 9498.571
$ws.Range("L43").Value = 9498.571
$ws.Range("N43").Value = -9636.571
$ws.Range("H74").Value = 20428304
$ws.Range("I74").Value = 23816356
$ws.Range("J74").Value = 100000
$ws.Range("K74").Value = 23816356
$ws.Range("L74").Value = 100000
$ws.Range("M74").Value = -23815420
$ws.Range("N74").Value = -101872
$ws.Range("H77").Value = 20428304
$ws.Range("I77").Value = 23816356
$ws.Range("J77").Value = 100000
$ws.Range("K77").Value = 119081780
$ws.Range("L77").Value = 500000
$ws.Range("M77").Value = -119077100
$ws.Range("N77").Value = -509360
$ws.Range("H141").Value = 3383.55
$ws.Range("I141").Value = 3092.889
$ws.Range("J141").Value = 5999.5
$ws.Range("K141").Value = 9278.667000000001
$ws.Range("L141").Value = 17998.5
$ws.Range("M141").Value = -4098.667000000001
$ws.Range("N141").Value = -28358.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1660613.4
$ws.Range("I32").Value = 735.64703
$ws.Range("J32").Value = 6951474
$ws.Range("K32").Value = 735.64703
$ws.Range("L32").Value = 6951474
$ws.Range("M32").Value = -448.64703
$ws.Range("N32").Value = -6952048
$ws.Range("H61").Value = 6698.769
$ws.Range("I61").Value = 8155.7144
$ws.Range("K61").Value = 8155.7144
$ws.Range("M61").Value = -7943.7144
$ws.Range("H88").Value = 6088.8
$ws.Range("J88").Value = 7020.5
$ws.Range("L88").Value = 7020.5
$ws.Range("N88").Value = -7832.5
$ws.Range("H91").Value = 6088.8
$ws.Range("J91").Value = 7020.5
$ws.Range("L91").Value = 7020.5
$ws.Range("N91").Value = -9828.5
$ws.Range("H97").Value = 1549.3334
$ws.Range("I97").Value = 1399.4
$ws.Range("K97").Value = 1399.4
$ws.Range("M97").Value = -903.4000000000001
$ws.Range("H122").Value = 6683.7
$ws.Range("I122").Value = 6576.5
$ws.Range("K122").Value = 19729.5
$ws.Range("M122").Value = -17279.5
$ws.Range("H135").Value = 44000
$ws.Range("J135").Value = 53500
$ws.Range("L135").Value = 53500
$ws.Range("N135").Value = -63640
$ws.Range("H136").Value = 6698.769
$ws.Range("I136").Value = 8155.7144
$ws.Range("K136").Value = 24467.1432
$ws.Range("M136").Value = -21917.1432

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1868.0769
$ws.Range("I22").Value = 1988.5555
$ws.Range("J22").Value = 1597
$ws.Range("K22").Value = 1988.5555
$ws.Range("L22").Value = 1597
$ws.Range("M22").Value = -1815.5555
$ws.Range("N22").Value = -1943
$ws.Range("H86").Value = 26641.2
$ws.Range("I86").Value = 34485.332
$ws.Range("J86").Value = 14875
$ws.Range("K86").Value = 34485.332
$ws.Range("L86").Value = 14875
$ws.Range("M86").Value = -33362.332
$ws.Range("N86").Value = -17121
$ws.Range("H89").Value = 26641.2
$ws.Range("I89").Value = 34485.332
$ws.Range("J89").Value = 14875
$ws.Range("K89").Value = 172426.66
$ws.Range("L89").Value = 74375
$ws.Range("M89").Value = -166810.66
$ws.Range("N89").Value = -85607
$ws.Range("H99").Value = 7688.8335
$ws.Range("I99").Value = 1955
$ws.Range("K99").Value = 1955
$ws.Range("M99").Value = -457

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 490651.47
$ws.Range("I22").Value = 943542
$ws.Range("K22").Value = 943542
$ws.Range("M22").Value = -943192
$ws.Range("H58").Value = 62516064
$ws.Range("I58").Value = 111122730
$ws.Range("J58").Value = 21782.857
$ws.Range("K58").Value = 111122730
$ws.Range("L58").Value = 21782.857
$ws.Range("M58").Value = -111122527
$ws.Range("N58").Value = -22188.857
$ws.Range("H99").Value = 5855643.5
$ws.Range("I99").Value = 11115148
$ws.Range("J99").Value = 11749.444
$ws.Range("K99").Value = 11115148
$ws.Range("L99").Value = 11749.444
$ws.Range("M99").Value = -11113650
$ws.Range("N99").Value = -14745.444
$ws.Range("H105").Value = 66673304
$ws.Range("I105").Value = 83339220
$ws.Range("J105").Value = 9666.666999999999
$ws.Range("K105").Value = 83339220
$ws.Range("L105").Value = 9666.666999999999
$ws.Range("M105").Value = -83337473
$ws.Range("N105").Value = -13160.667
$ws.Range("H122").Value = 5580.5264
$ws.Range("I122").Value = 5284.6
$ws.Range("J122").Value = 5909.3335
$ws.Range("K122").Value = 15853.8
$ws.Range("L122").Value = 17728.0005
$ws.Range("M122").Value = -13403.8
$ws.Range("N122").Value = -22628.0005
$ws.Range("H126").Value = 5855643.5
$ws.Range("I126").Value = 11115148
$ws.Range("J126").Value = 11749.444
$ws.Range("K126").Value = 33345444
$ws.Range("L126").Value = 35248.33199999999
$ws.Range("M126").Value = -33342974
$ws.Range("N126").Value = -40188.33199999999
$ws.Range("H136").Value = 62516064
$ws.Range("I136").Value = 111122730
$ws.Range("J136").Value = 21782.857
$ws.Range("K136").Value = 333368190
$ws.Range("L136").Value = 65348.571
$ws.Range("M136").Value = -333365640
$ws.Range("N136").Value = -70448.571

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 249998
$ws.Range("J37").Value = 249998
$ws.Range("L37").Value = 749994
$ws.Range("N37").Value = -750218
$ws.Range("H56").Value = 10699.667
$ws.Range("I56").Value = 10699.667
$ws.Range("K56").Value = 10699.667
$ws.Range("M56").Value = -10169.667
$ws.Range("H58").Value = 6999.6665
$ws.Range("I58").Value = 999
$ws.Range("K58").Value = 2997
$ws.Range("M58").Value = -2869
$ws.Range("H132").Value = 1896.7742
$ws.Range("I132").Value = 894.4706
$ws.Range("K132").Value = 8050.2354
$ws.Range("M132").Value = -5520.2354

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 75666.664
$ws.Range("J125").Value = 75666.664
$ws.Range("L125").Value = 75666.664
$ws.Range("N125").Value = -80586.664

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6386.8887
$ws.Range("I40").Value = 4834.3335
$ws.Range("J40").Value = 7163.1665
$ws.Range("K40").Value = 4834.3335
$ws.Range("L40").Value = 7163.1665
$ws.Range("M40").Value = -4698.3335
$ws.Range("N40").Value = -7435.1665
$ws.Range("H46").Value = 38462550
$ws.Range("I46").Value = 1166.6666
$ws.Range("J46").Value = 50000970
$ws.Range("K46").Value = 1166.6666
$ws.Range("L46").Value = 50000970
$ws.Range("M46").Value = -978.6666
$ws.Range("N46").Value = -50001346
$ws.Range("H61").Value = 4249.1143
$ws.Range("I61").Value = 2899.0334
$ws.Range("J61").Value = 12349.6
$ws.Range("K61").Value = 2899.0334
$ws.Range("L61").Value = 12349.6
$ws.Range("M61").Value = -2697.0334
$ws.Range("N61").Value = -12753.6
$ws.Range("H68").Value = 2683.2778
$ws.Range("I68").Value = 2393.6875
$ws.Range("K68").Value = 2393.6875
$ws.Range("M68").Value = -1644.6875
$ws.Range("H71").Value = 2683.2778
$ws.Range("I71").Value = 2393.6875
$ws.Range("K71").Value = 11968.4375
$ws.Range("M71").Value = -8224.4375
$ws.Range("H113").Value = 4249.1143
$ws.Range("I113").Value = 2899.0334
$ws.Range("J113").Value = 12349.6
$ws.Range("K113").Value = 2899.0334
$ws.Range("L113").Value = 12349.6
$ws.Range("M113").Value = -729.0333999999998
$ws.Range("N113").Value = -16689.6
$ws.Range("H122").Value = 9472
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 9472
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 28416
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -33316

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 44965.668
$ws.Range("J75").Value = 44965.668
$ws.Range("L75").Value = 44965.668
$ws.Range("N75").Value = -46837.668
$ws.Range("H78").Value = 44965.668
$ws.Range("J78").Value = 44965.668
$ws.Range("L78").Value = 134897.004
$ws.Range("N78").Value = -144257.004
$ws.Range("H100").Value = 927.2727
$ws.Range("I100").Value = 671.5714
$ws.Range("K100").Value = 1343.1428
$ws.Range("M100").Value = -802.1428000000001
$ws.Range("H104").Value = 21471.9
$ws.Range("J104").Value = 21471.9
$ws.Range("L104").Value = 21471.9
$ws.Range("N104").Value = -28459.9
$ws.Range("H122").Value = 3793.6667
$ws.Range("I122").Value = 3892.875
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 11678.625
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -9228.625
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 6375.6
$ws.Range("I126").Value = 5630.9165
$ws.Range("J126").Value = 9354.333000000001
$ws.Range("K126").Value = 16892.7495
$ws.Range("L126").Value = 28062.999
$ws.Range("M126").Value = -14422.7495
$ws.Range("N126").Value = -33002.999
$ws.Range("H136").Value = 23815374
$ws.Range("I136").Value = 26316742
$ws.Range("K136").Value = 78950226
$ws.Range("M136").Value = -78947676
$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200
